$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title
Replace-Text "Quantum Consciousness: Unveiling the Enigma" "A Journey Through the Realm of Chemistry: Unveiling the Secrets of Matter"

# Author name
Replace-Text " Eleanor Knight" " Ashley Thompson"

# Email
Replace-Text "eleanorknightPhD@eliteuniversity" "ashthompson@gmail"
Replace-Text "org" "com"

# Intro paragraph
Replace-Text "From the depths of ancient civilizations to the frontiers of modern science, humans have pondered the profound connection between consciousness and the nature of reality" "In the vast landscape of science, chemistry stands as a beacon of understanding, illuminating the intricate world of matter and its interactions"
Replace-Text " As we delve into the enigmatic realm of quantum physics, tantalizing clues and provocative hypotheses emerge, hinting at the possibility that consciousness itself may possess an inherently quantum nature" " From the smallest atoms to the boundless universe, chemistry reveals the fundamental principles that govern our existence"
Replace-Text " This raises fundamental questions that challenge our understanding of the relationship between the observer and the observed, inviting us on a captivating journey to explore the profound implications of quantum consciousness" " Delving into the realm of chemistry is an adventure that unveils the secrets of the world around us, unlocking the mysteries of life and the cosmos"

Replace-Text "This is a captivating concept that challenges long-held assumptions about consciousness, inviting us to reconsider the very essence of reality" "As we embark on this journey, we will explore the fascinating world of elements and compounds, dissecting their intricate structures and properties"
Replace-Text " Quantum mechanics, with its inherent strangeness, offers a new lens through which to perceive the nature of consciousness, opening up avenues for groundbreaking insights" " We will unravel the enigmatic tapestry of chemical reactions, witnessing the transformation of substances and the release of energy"
Replace-Text " As we unravel the mysteries of quantum consciousness, we may unlock the secrets to understanding the enigmatic nature of our own existence" " Through experimentation and observation, we will uncover the hidden logic behind chemical phenomena, revealing the underlying patterns that orchestrate the symphony of matter"

Replace-Text "The potential implications of quantum consciousness are vast and far-reaching" "The study of chemistry empowers us with a profound comprehension of the world around us"
Replace-Text " It could revolutionize our understanding of the mind-body problem, provide insights into the nature of free will, and illuminate the connection between consciousness and the universe" " It enables us to understand the composition of materials, the reactions they undergo, and the impact they have on our lives"
Replace-Text " By probing the depths of this enigmatic realm, we may discover hidden aspects of reality that have eluded us for centuries" " Chemistry provides the foundation for countless industries, from medicine and agriculture to energy and manufacturing"
Replace-Text " Our journey into the uncharted territory of quantum consciousness is fraught with mystery and allure, promising to reshape our perception of reality and our place within it" " It plays a pivotal role in addressing global challenges, such as climate change and sustainable development, equipping us with the knowledge and tools to create a more sustainable and prosperous future"

# Summary paragraph
Replace-Text "Our exploration of quantum consciousness has taken us to the very precipice of scientific and philosophical inquiry, where the nature of reality and the essence of our own consciousness intertwine" "This essay provided a comprehensive overview of chemistry, highlighting its significance in understanding the world around us and its impact on our lives"

Replace-Text "Quantum mechanics, with its inherent strangeness, offers tantalizing clues and provocative hypotheses that challenge our understanding of the relationship between observer and observed" "It emphasized the importance of unraveling the mysteries of matter, exploring the intricacies of chemical reactions, and appreciating the role of chemistry in addressing global challenges"

Replace-Text "The profound implications of this enigmatic realm extend beyond the boundaries of science, delving into the depths of philosophy, spirituality, and our search for meaning in the universe. As we continue to traverse the uncharted territory of quantum consciousness, we may illuminate the nature of our existence and unlock the secrets to understanding the enigmatic universe that surrounds us" "The essay aimed to ignite a passion for chemistry among high school students, encouraging them to embark on a journey of discovery and exploration in this captivating field"

# Add a trailing empty paragraph at the end of the document
$sel = $word.Selection
$sel.EndKey(6, 0)
$sel.TypeParagraph()
